$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notifications")

# Row 4: On Second Vote - mark Status/Api as done
$ws.Range("E4").Value = "Done"
$ws.Range("F4").Value = "Api"

# Row 10: On Create Challage - update Send To and Link, mark as done / Api / Web
$ws.Range("B10").Value = "to Challenge Owner, to Admin"
$ws.Range("C10").Value = "CHALLENGE_DETAIL_SCREEN"
$ws.Range("E10").Value = "Done"
$ws.Range("F10").Value = "Api"
$ws.Range("G10").Value = "Web"

# Scroll the sheet view so row 2 is at the top
$ws.Application.ActiveWindow.ScrollRow = 2
